$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '28.335.75'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +4.08%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.732.48'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +2.62%  '

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.13%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '219.46'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +1.56%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.523'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +0.20%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.14%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '24.15'

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.268'

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0637'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +1.52%  '

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +0.54%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.975.26'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +2.56%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.721.18'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +2.00%  '

$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +1.58%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.564'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +2.15%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '67.83'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +0.67%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '28.318.64'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +4.09%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '247.68'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +4.02%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.0₃0755'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +1.27%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.93'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -2.94%  '

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -0.12%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.66'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +1.66%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.70'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +0.22%  '

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -0.73%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '149.34'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +0.63%  '

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +2.54%  '

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +1.30%  '

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +0.47%  '

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -0.31%  '

$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +2.67%  '

$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +2.62%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.43'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +0.73%  '

$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +0.60%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.490.25'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -5.84%  '

$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -2.21%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.979'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +1.94%  '

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -0.32%  '

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +0.59%  '

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +1.13%  '

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +0.44%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '70.32'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +0.96%  '

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -0.14%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.66'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -2.44%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.879.70'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +2.35%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.30'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +1.43%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.798'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +1.31%  '

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +7.22%  '

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +3.75%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '90.52'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -0.97%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '8.19'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -0.61%  '

$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -0.83%  '
